$d = $word.ActiveDocument

# The paragraph ending "...merece la vuestra grandeza»." is the last piece
# of body text before the document's trailing empty paragraph. We insert a
# brand-new paragraph right after it (inheriting that paragraph's Courier
# New formatting automatically), containing the text "djioasfjiojhioas".

$anchor = $d.Content
$found = $anchor.Find.Execute("merece la vuestra grandeza».", $true, $false, `
    $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $anchorParaIndex = $anchor.Paragraphs(1).Index
} else {
    # Fallback: anchor on the last non-empty paragraph.
    $anchorParaIndex = $d.Paragraphs.Count - 1
}

$insertAt = $d.Paragraphs($anchorParaIndex).Range
$insertAt.Collapse(0)                 # wdCollapseEnd = 0
$insertAt.InsertParagraphAfter() | Out-Null

$newPara = $d.Paragraphs($anchorParaIndex + 1).Range
$newPara.Text = "djioasfjiojhioas"
